# Update the "Förändrad" (Changed) date column (C) for all data rows
# from 2023-09-09 (45178) to 2023-09-10 (45179).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 200; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45178) {
        $cell.Value2 = 45179
    }
}
